# Adding/moving preprocessed files to "preprocessed" folder.
# This script adds two new columns to the regions table:
#   G: "refseq"   - the RefSeq NC_ accession for the chromosome of each region
#   H: "crispick" - a CRISPick-style locus string "<accession>:+:<start>-<end>"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map chromosome name (as it appears in column A) to its RefSeq NC_ accession.
$chrMap = @{
    "11" = "NC_000011.9"
    "8"  = "NC_000008.10"
    "16" = "NC_000016.9"
    "6"  = "NC_000006.11"
    "17" = "NC_000017.10"
    "19" = "NC_000019.9"
    "7"  = "NC_000007.13"
    "12" = "NC_000012.11"
    "10" = "NC_000010.10"
    "9"  = "NC_000009.11"
    "20" = "NC_000020.10"
    "14" = "NC_000014.8"
    "2"  = "NC_000002.11"
    "5"  = "NC_000005.9"
    "15" = "NC_000015.9"
    "21" = "NC_000021.8"
    "4"  = "NC_000004.11"
}

# Header row
$ws.Cells.Item(1, 7).Value = "refseq"
$ws.Cells.Item(1, 8).Value = "crispick"

# Find the last populated row in column A (chr column)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $chr = [string]$ws.Cells.Item($r, 1).Value()
    $start = $ws.Cells.Item($r, 2).Value()
    $end = $ws.Cells.Item($r, 3).Value()

    $accession = $chrMap[$chr]

    $ws.Cells.Item($r, 7).Value = $accession
    $ws.Cells.Item($r, 8).Value = "$($accession):+:$([int64]$start)-$([int64]$end)"
}
